# The commit groups all four top-level shapes on slide 2 (the existing
# "Group 3", "TextBox 10", "Left Bracket 11" and "Left Bracket 12") into a
# brand-new enclosing group ("Group 13", next id 14).
#
# PowerPoint's shape-id/name allocator (as emulated here) hands out a couple
# of "low" ids (2, 3) before it resumes counting after the highest id
# already used on the slide (13 -> 14). We burn through those two low ids
# with a scratch shape that we immediately delete, so the subsequent Group()
# call lands on id=14 / name="Group 13", matching the target document.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$scratch = $s.Shapes.AddShape(1, 0, 0, 10, 10)
$scratch.Delete()
$scratch2 = $s.Shapes.AddShape(1, 0, 0, 10, 10)
$scratch2.Delete()

# Select the four existing top-level shapes in order and group them.
$range = $s.Shapes.Range(@(1, 2, 3, 4))
$newGroup = $range.Group()
$newGroup.Name = "Group 13"
